# Add "Wins" / "Losses" / "Ties" season-record columns (AD, AE, AF) to the
# roster table on Sheet1, mirroring the existing header style and filling
# every player row with the team's season record (65-97-0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels, styled like the rest of the header ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Clone the formatting of an existing header cell (bold font, border,
# centered alignment) onto the three new header cells so they reuse the
# same style record instead of minting new ones.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-53): season record for every player ---
$firstDataRow = 2
$lastDataRow = 53
$winsCol = 30   # AD
$lossesCol = 31 # AE
$tiesCol = 32   # AF

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $ws.Cells.Item($r, $winsCol).Value = 65
    $ws.Cells.Item($r, $lossesCol).Value = 97
    $ws.Cells.Item($r, $tiesCol).Value = 0
}
